$wb = $excel.ActiveWorkbook

# --- CmsWork sheet ---
$ws = $wb.Worksheets.Item("CmsWork")

# Row 2 (http://example.com/collection0/work1)
$ws.Range("C2").Value = "_:N3e06727f8e624edfae132632f64eca90"
$ws.Range("D2").Value = "http://example.com/organization3"

# Row 3 (http://example.com/collection0/work3)
$ws.Range("C3").Value = "_:Neb9ed50049964b01aa889cafbaebca74"
$ws.Range("U3").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:103"

# Row 4 (http://example.com/collection1/work5)
$ws.Range("C4").Value = "_:N94d92baeef54447b973a44af77cd60ee"
$ws.Range("D4").Value = "http://example.com/person2"
$ws.Range("E4").Value = "CmsCollection1CmsWork5 alternative title 0"

# Row 5 (http://example.com/collection1/work7)
$ws.Range("C5").Value = "_:N789d7df36c734532a8c7b501098c00d4"
$ws.Range("D5").Value = "http://example.com/organization0"
$ws.Range("E5").Value = "CmsCollection1CmsWork7 alternative title 1"
$ws.Range("J5").Value = "CmsCollection1CmsWork7 provenance 1"

# Row 6 (http://example.com/freestandingwork9)
$ws.Range("B6").Value = "_:N7a67edd2c6a14d6e8a4c1e4d977cd670"
$ws.Range("C6").Value = "http://example.com/organization1"
$ws.Range("D6").Value = "FreestandingWork9 alternative title 1"
$ws.Range("I6").Value = "FreestandingWork9 provenance 0"
$ws.Range("T6").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:100"

# Row 7 (http://example.com/freestandingwork11)
$ws.Range("B7").Value = "_:Na19116e1d54c4b6d91eaa893ebbd0162"
$ws.Range("C7").Value = "http://example.com/organization3"
$ws.Range("D7").Value = "FreestandingWork11 alternative title 1"
$ws.Range("I7").Value = "FreestandingWork11 provenance 1"

# --- CmsWorkClosing sheet ---
$wsClosing = $wb.Worksheets.Item("CmsWorkClosing")

$wsClosing.Range("A2").Value = "_:N00f42f6db6b242d3aefb0ee047ddc7ae"
$wsClosing.Range("C2").Value = "_:N171f1109b61b43fa81f161f6fe15db00"

$wsClosing.Range("A3").Value = "_:Ncaf0a51df02146f891c2aec439b1ea9b"
$wsClosing.Range("C3").Value = "_:N81647bd550c94976bad78b3afbcbb366"

$wsClosing.Range("A4").Value = "_:Nfbec17ceba9d4e27b1f5200a39ca4730"
$wsClosing.Range("C4").Value = "_:Na67ce142815f4520be2c18b0ca17b92f"

$wsClosing.Range("A5").Value = "_:Ncf152b4be0df4e70953e5705c667c5aa"
$wsClosing.Range("C5").Value = "_:N20ba0e867e5f4e2a9168aba47abaa781"

$wsClosing.Range("A6").Value = "_:Ncf27c1254e804042b374217a72661e62"
$wsClosing.Range("C6").Value = "_:N44b6abbf540e4120afc737eba9d61cae"

$wsClosing.Range("A7").Value = "_:N7702e8ab0bd04acf83c63c001821b0fb"
$wsClosing.Range("C7").Value = "_:N57daa86fd1b34106825f3af81e3d3c65"

# --- CmsWorkOpening sheet ---
$wsOpening = $wb.Worksheets.Item("CmsWorkOpening")

$wsOpening.Range("C2").Value = "_:N171f1109b61b43fa81f161f6fe15db00"
$wsOpening.Range("C3").Value = "_:N81647bd550c94976bad78b3afbcbb366"
$wsOpening.Range("C4").Value = "_:Na67ce142815f4520be2c18b0ca17b92f"
$wsOpening.Range("C5").Value = "_:N20ba0e867e5f4e2a9168aba47abaa781"
$wsOpening.Range("C6").Value = "_:N44b6abbf540e4120afc737eba9d61cae"
$wsOpening.Range("C7").Value = "_:N57daa86fd1b34106825f3af81e3d3c65"

Write-Output "done"
